$d = $word.ActiveDocument

# Pull the whole package as flat OOXML so we can surgically rewrite the
# <w:docDefaults> block inside word/styles.xml (there is no dedicated
# Styles.DocDefaults object in the Word OM).
$xml = $d.WordOpenXML

$pattern = "(?s)<w:docDefaults>.*?</w:docDefaults>"

$replacement = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newXml = $xml -replace $pattern, $replacement

$d.WordOpenXML = $newXml

Write-Output "docDefaults rewritten"
